# Refresh the cryptocurrency table: updates Price (column D) and
# Volume(1h) (column E) for every coin, and reflects EnergySwap overtaking
# Decentraland in rows 50-51 (name, link, price and volume all swap rank).
# Matches the "Updated cryptos list ... with GitHub Actions" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '25.907.88' },
    @{ Cell = 'E2'; Value = '  -0.71%  ' },
    @{ Cell = 'D3'; Value = '1.741.34' },
    @{ Cell = 'E3'; Value = '  -0.47%  ' },
    @{ Cell = 'D4'; Value = '1.000' },
    @{ Cell = 'E4'; Value = '  +0.03%  ' },
    @{ Cell = 'D5'; Value = '248.29' },
    @{ Cell = 'E5'; Value = '  +5.65%  ' },
    @{ Cell = 'D6'; Value = '1.000' },
    @{ Cell = 'E6'; Value = '  +0.03%  ' },
    @{ Cell = 'D7'; Value = '0.5083' },
    @{ Cell = 'E7'; Value = '  -3.75%  ' },
    @{ Cell = 'D8'; Value = '0.2719' },
    @{ Cell = 'E8'; Value = '  -2.58%  ' },
    @{ Cell = 'D9'; Value = '0.06180' },
    @{ Cell = 'E9'; Value = '  -0.18%  ' },
    @{ Cell = 'D10'; Value = '1.743.44' },
    @{ Cell = 'E10'; Value = '  -0.29%  ' },
    @{ Cell = 'D11'; Value = '0.07234' },
    @{ Cell = 'E11'; Value = '  +0.72%  ' },
    @{ Cell = 'D12'; Value = '15.12' },
    @{ Cell = 'E12'; Value = '  -1.86%  ' },
    @{ Cell = 'D13'; Value = '0.6476' },
    @{ Cell = 'E13'; Value = '  +0.29%  ' },
    @{ Cell = 'D14'; Value = '4.623' },
    @{ Cell = 'E14'; Value = '  -0.07%  ' },
    @{ Cell = 'D15'; Value = '77.59' },
    @{ Cell = 'E15'; Value = '  -1.15%  ' },
    @{ Cell = 'E16'; Value = '  +0.06%  ' },
    @{ Cell = 'D17'; Value = '1.000' },
    @{ Cell = 'E17'; Value = '  +0.02%  ' },
    @{ Cell = 'D18'; Value = '25.931.96' },
    @{ Cell = 'E18'; Value = '  -0.24%  ' },
    @{ Cell = 'D20'; Value = '0.000006808' },
    @{ Cell = 'E20'; Value = '  +1.20%  ' },
    @{ Cell = 'D21'; Value = '1.972.24' },
    @{ Cell = 'E21'; Value = '  +0.18%  ' },
    @{ Cell = 'D22'; Value = '4.281' },
    @{ Cell = 'E22'; Value = '  -0.72%  ' },
    @{ Cell = 'E23'; Value = '  -1.12%  ' },
    @{ Cell = 'D24'; Value = '5.391' },
    @{ Cell = 'E24'; Value = '  +2.90%  ' },
    @{ Cell = 'D25'; Value = '136.26' },
    @{ Cell = 'E25'; Value = '  -1.20%  ' },
    @{ Cell = 'E26'; Value = '  -0.65%  ' },
    @{ Cell = 'E27'; Value = '  -0.41%  ' },
    @{ Cell = 'D28'; Value = '1.774' },
    @{ Cell = 'E28'; Value = '  -1.65%  ' },
    @{ Cell = 'D29'; Value = '105.46' },
    @{ Cell = 'E29'; Value = '  +0.94%  ' },
    @{ Cell = 'D30'; Value = '3.909' },
    @{ Cell = 'E30'; Value = '  +2.91%  ' },
    @{ Cell = 'D31'; Value = '0.08231' },
    @{ Cell = 'E31'; Value = '  -0.61%  ' },
    @{ Cell = 'D32'; Value = '3.636' },
    @{ Cell = 'E32'; Value = '  -0.50%  ' },
    @{ Cell = 'E33'; Value = '  +2.47%  ' },
    @{ Cell = 'D34'; Value = '2.656' },
    @{ Cell = 'E34'; Value = '  +0.41%  ' },
    @{ Cell = 'D35'; Value = '0.9951' },
    @{ Cell = 'E35'; Value = '  -1.02%  ' },
    @{ Cell = 'D36'; Value = '0.6240' },
    @{ Cell = 'E36'; Value = '  -1.71%  ' },
    @{ Cell = 'D37'; Value = '2.731' },
    @{ Cell = 'E37'; Value = '  +0.76%  ' },
    @{ Cell = 'D38'; Value = '0.01601' },
    @{ Cell = 'E38'; Value = '  +0.21%  ' },
    @{ Cell = 'D39'; Value = '1.919' },
    @{ Cell = 'E39'; Value = '  -1.45%  ' },
    @{ Cell = 'D40'; Value = '0.9994' },
    @{ Cell = 'E40'; Value = '  +0.00%  ' },
    @{ Cell = 'D41'; Value = '99.30' },
    @{ Cell = 'E41'; Value = '  -1.27%  ' },
    @{ Cell = 'D42'; Value = '0.7580' },
    @{ Cell = 'E42'; Value = '  +1.73%  ' },
    @{ Cell = 'D43'; Value = '0.3848' },
    @{ Cell = 'E43'; Value = '  -2.01%  ' },
    @{ Cell = 'D44'; Value = '4.996' },
    @{ Cell = 'E44'; Value = '  -0.58%  ' },
    @{ Cell = 'D45'; Value = '0.1134' },
    @{ Cell = 'E45'; Value = '  -1.04%  ' },
    @{ Cell = 'D46'; Value = '6.287' },
    @{ Cell = 'E46'; Value = '  -0.99%  ' },
    @{ Cell = 'D47'; Value = '55.51' },
    @{ Cell = 'E47'; Value = '  +2.36%  ' },
    @{ Cell = 'D48'; Value = '0.05237' },
    @{ Cell = 'E48'; Value = '  -2.10%  ' },
    @{ Cell = 'E49'; Value = '  -0.77%  ' },
    @{ Cell = 'B50'; Value = 'EnergySwap' },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' },
    @{ Cell = 'D50'; Value = '7.514' },
    @{ Cell = 'E50'; Value = '  -1.33%  ' },
    @{ Cell = 'B51'; Value = 'Decentraland' },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana' },
    @{ Cell = 'D51'; Value = '7.474' },
    @{ Cell = 'E51'; Value = '  -1.14%  ' }
)

foreach ($update in $updates) {
    $range = $ws.Range($update.Cell)
    $value = $update.Value

    # Columns D/E were authored as literal text ("1.000", "0.06180",
    # "  -0.71%  ", "25.907.88", ...) rather than numbers, so exact
    # formatting (trailing zeros, padding, stray "thousand-dot" typos in
    # the source feed, etc.) has to survive the write. Assigning a bare
    # numeric-looking string to .Value gets auto-coerced to a Number by
    # Excel, silently dropping that formatting - so for anything that is
    # purely digits/one-decimal-point we force text via the classic
    # leading-apostrophe trick and then restore the cell's original
    # (default) style so no stray quote-prefix formatting is left behind.
    if ($value -match '^[0-9]+(\.[0-9]+)?$') {
        $range.Value = "'" + $value
        $range.Style = 'Normal'
    } else {
        $range.Value = $value
    }
}
